$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 5231.579170919113
$ws.Range("C2").Value = 5608.494317741896
$ws.Range("B3").Value = 5227.18998025917
$ws.Range("C3").Value = 5598.945408379392
$ws.Range("B4").Value = 5226.135994334955
$ws.Range("C4").Value = 5595.360088523847
$ws.Range("B5").Value = 5218.753282791317
$ws.Range("C5").Value = 5584.401167647237
$ws.Range("B6").Value = 5213.326815717922
$ws.Range("C6").Value = 5586.692958707069
$ws.Range("B7").Value = 5181.597320037057
$ws.Range("C7").Value = 5526.517831907375
$ws.Range("B8").Value = 5148.424771118935
$ws.Range("C8").Value = 5489.186722942089
$ws.Range("B9").Value = 5117.056795487993
$ws.Range("C9").Value = 5462.248588575165
$ws.Range("B10").Value = 4959.100050252483
$ws.Range("C10").Value = 5310.349150550828
$ws.Range("B11").Value = 4897.697391440604
$ws.Range("C11").Value = 5256.034139637408
$ws.Range("B12").Value = 4877.112749153162
$ws.Range("C12").Value = 5236.263592709773
$ws.Range("B13").Value = 4870.357160500368
$ws.Range("C13").Value = 5243.170989172474
$ws.Range("B14").Value = 4854.419310366724
$ws.Range("C14").Value = 5234.243923314793
$ws.Range("B15").Value = 4843.045826983373
$ws.Range("C15").Value = 5212.387337925913
$ws.Range("B16").Value = 4843.853480301618
$ws.Range("C16").Value = 5214.068932723902
$ws.Range("B17").Value = 4844.941348929079
$ws.Range("C17").Value = 5219.726437521515
$ws.Range("B18").Value = 4825.520421821459
$ws.Range("C18").Value = 5178.99899342375
$ws.Range("B19").Value = 4803.236894393954
$ws.Range("C19").Value = 5174.888737863336
$ws.Range("B20").Value = 4802.579039721439
$ws.Range("C20").Value = 5172.126452361473
$ws.Range("B21").Value = 4794.230050057236
$ws.Range("C21").Value = 5144.745724301809
$ws.Range("B22").Value = 4785.892576676623
$ws.Range("C22").Value = 5147.564658529823
$ws.Range("B23").Value = 4784.476059949327
$ws.Range("C23").Value = 5147.80473789105
$ws.Range("B24").Value = 4784.146700098136
$ws.Range("C24").Value = 5144.722616980358
$ws.Range("B25").Value = 4763.615903825849
$ws.Range("C25").Value = 5125.158938291903
$ws.Range("B26").Value = 4763.603567251298
$ws.Range("C26").Value = 5124.650345156093
$ws.Range("B27").Value = 4763.603758122062
$ws.Range("C27").Value = 5124.643557621705
